$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 317, pushing the existing rows 317:324 down to 318:325
$ws.Rows.Item(317).Insert()

# Populate the newly inserted row 317 with the new weekly price record
$ws.Cells.Item(317, 1).Value = 9
$ws.Cells.Item(317, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(317, 3).Value = "Metropolitana"
$ws.Cells.Item(317, 4).Value = 45075
$ws.Cells.Item(317, 5).Value = 13
$ws.Cells.Item(317, 6).Value = 100112003
$ws.Cells.Item(317, 7).Value = "Ajo"
$ws.Cells.Item(317, 8).Value = "Chino"
$ws.Cells.Item(317, 9).Value = "Primera"
$ws.Cells.Item(317, 10).Value = 610
$ws.Cells.Item(317, 11).Value = 14000
$ws.Cells.Item(317, 12).Value = 15000
$ws.Cells.Item(317, 13).Value = 14500
$ws.Cells.Item(317, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(317, 15).Value = "China"
$ws.Cells.Item(317, 16).Value = 1450
$ws.Cells.Item(317, 17).Value = 10
$ws.Cells.Item(317, 18).Value = "Hortaliza"
